$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Septiembre de 2020 a las 05:14"

# Update country data rows (B=Casos totales, C=Nuevos casos, D=Casos activos,
# E=Recuperados, F=Casos criticos, G=Muertes hoy, H=Muertes)

# Row 8: Peru
$ws.Range("B8").Value = 676848
$ws.Range("D8").Value = 498523
$ws.Range("E8").Value = 148771
$ws.Range("H8").Value = 29554

# Row 31: Bolivia
$ws.Range("B31").Value = 119580
$ws.Range("C31").Value = 799
$ws.Range("D31").Value = 67844
$ws.Range("E31").Value = 46393
$ws.Range("G31").Value = 55
$ws.Range("H31").Value = 5343

# Row 33: Kazajistan
$ws.Range("B33").Value = 106225
$ws.Range("C33").Value = 104
$ws.Range("E33").Value = 6052

# Row 39: Belgica
$ws.Range("B39").Value = 87174
$ws.Range("C39").Value = 630
$ws.Range("D39").Value = 18534
$ws.Range("E39").Value = 58739
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 9901

# Row 50: Honduras
$ws.Range("B50").Value = 63798
$ws.Range("C50").Value = 640
$ws.Range("D50").Value = 12347
$ws.Range("E50").Value = 49467
$ws.Range("G50").Value = 30
$ws.Range("H50").Value = 1984

# Row 73: Australia
$ws.Range("B73").Value = 26207
$ws.Range("C73").Value = 71
$ws.Range("D73").Value = 22172
$ws.Range("E73").Value = 3287
$ws.Range("G73").Value = 11
$ws.Range("H73").Value = 748

# Row 133: Jamaica
$ws.Range("B133").Value = 2964
$ws.Range("C133").Value = 68
$ws.Range("D133").Value = 922
$ws.Range("E133").Value = 2012
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 30

# Row 172: Islas Turcas y Caicos
$ws.Range("B172").Value = 591
$ws.Range("C172").Value = 14
$ws.Range("D172").Value = 221
$ws.Range("E172").Value = 365
